# Generate Report for Handback
# Adds a new handback entry (862cc7d4-a310-46be-aa6b-6233911e0e15) as row 4
# to the Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview" (table3) - columns: A File Name, B Path And Name,
# C Extension, D Publish URL, E zh-cn, F de-de, G Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value2 = "862cc7d4-a310-46be-aa6b-6233911e0e15.md"
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/403b1ea504cb66a770324e8b54e8b78909f303c3/e2e/862cc7d4-a310-46be-aa6b-6233911e0e15.md",
    $null,
    $null,
    "e2e\862cc7d4-a310-46be-aa6b-6233911e0e15.md"
) | Out-Null
$wsOverview.Range("C4").Value2 = ".md"
$wsOverview.Range("E4").Value2 = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value2 = "Handed back: in sync with en-US"
$wsOverview.Range("G4").Value2 = "2016-11-08 22:41:55"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (table1) - columns:
# A Source File Name, B File Extension, C Status, D Source Path, E Priority,
# F Content Duplicate, G Correspond Handoff File, H Correspond Handoff Datetime,
# I Target File, J Correspond Handback File, K Correspond Handback DateTime,
# L Reference Tokens, M To be localized, N Dependency From, O Has metadata,
# P Error Detail
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/403b1ea504cb66a770324e8b54e8b78909f303c3/e2e/862cc7d4-a310-46be-aa6b-6233911e0e15.md",
    $null,
    $null,
    "862cc7d4-a310-46be-aa6b-6233911e0e15.md"
) | Out-Null
$wsZhCn.Range("B4").Value2 = ".md"
$wsZhCn.Range("C4").Value2 = "Handed back: in sync with en-US"
$wsZhCn.Range("D4").Value2 = "e2e"
$wsZhCn.Range("E4").Value2 = "ht"
$wsZhCn.Range("F4").Value2 = "True"
$wsZhCn.Range("G4").Value2 = "862cc7d4-a310-46be-aa6b-6233911e0e15.403b1ea504cb66a770324e8b54e8b78909f303c3.zh-cn.xlf"
$wsZhCn.Range("H4").Value2 = "2016-11-08 22:41:41"
$wsZhCn.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/403b1ea504cb66a770324e8b54e8b78909f303c3/e2e/862cc7d4-a310-46be-aa6b-6233911e0e15.md",
    $null,
    $null,
    "862cc7d4-a310-46be-aa6b-6233911e0e15.md"
) | Out-Null
$wsZhCn.Range("J4").Value2 = "862cc7d4-a310-46be-aa6b-6233911e0e15.403b1ea504cb66a770324e8b54e8b78909f303c3.zh-cn.xlf"
$wsZhCn.Range("K4").Value2 = "2016-11-08 22:42:34"
$wsZhCn.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("M4").Value2 = "True"
$wsZhCn.Range("O4").Value2 = "False"

# ---------------------------------------------------------------------------
# Sheet "de-de" (table2) - same column layout as zh-cn
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/403b1ea504cb66a770324e8b54e8b78909f303c3/e2e/862cc7d4-a310-46be-aa6b-6233911e0e15.md",
    $null,
    $null,
    "862cc7d4-a310-46be-aa6b-6233911e0e15.md"
) | Out-Null
$wsDeDe.Range("B4").Value2 = ".md"
$wsDeDe.Range("C4").Value2 = "Handed back: in sync with en-US"
$wsDeDe.Range("D4").Value2 = "e2e"
$wsDeDe.Range("E4").Value2 = "ht"
$wsDeDe.Range("F4").Value2 = "True"
$wsDeDe.Range("G4").Value2 = "862cc7d4-a310-46be-aa6b-6233911e0e15.403b1ea504cb66a770324e8b54e8b78909f303c3.de-de.xlf"
$wsDeDe.Range("H4").Value2 = "2016-11-08 22:41:55"
$wsDeDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/403b1ea504cb66a770324e8b54e8b78909f303c3/e2e/862cc7d4-a310-46be-aa6b-6233911e0e15.md",
    $null,
    $null,
    "862cc7d4-a310-46be-aa6b-6233911e0e15.md"
) | Out-Null
$wsDeDe.Range("J4").Value2 = "862cc7d4-a310-46be-aa6b-6233911e0e15.403b1ea504cb66a770324e8b54e8b78909f303c3.de-de.xlf"
$wsDeDe.Range("K4").Value2 = "2016-11-08 22:42:52"
$wsDeDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("M4").Value2 = "True"
$wsDeDe.Range("O4").Value2 = "False"
